$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openTickets")

$ws.Range("E2").Value = 'Analysis for Excel" add-in causes crashes.'
$ws.Range("E3").Value = "Cisco Softphone does not work"
$ws.Range("E4").Value = "Installing Nuance Software"
$ws.Range("E5").Value = "Thank you for your participation!"
$ws.Range("F5").Value = 'Thank you for supporting our research in artificial intelligence!"'
